$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D15").Value2 = "- Completed the introduction for the initial draft report"
$ws.Range("D20").Value2 = "N/A"
$ws.Range("E20").Value2 = "N/A"
$ws.Range("D30").Value2 = "- Look over colleted datasets once all have been compiled."
$ws.Range("E30").Value2 = "- Provide feedback on datasets and recommend which datasets use for further analysis"
$ws.Range("D35").Value2 = "-Reviewed Feedback given from Lecturers"
$ws.Range("E35").Value2 = "- Planned Next steps (When to do presentation & Cleaning datasets for EDA)"
$ws.Range("D45").Value2 = "-EDA / Cleaning a dataset as well as provide data visualisations for analysis`n- Discussed clustering model`n- Began Analysis of EDA (In relation to the draft project report for the client)"
$ws.Range("E45").Value2 = "- Carry on with EDA and cleaning`n- Get started on the presentation`n- Begin Analysis of EDA"
$ws.Range("D46").Value2 = "- Discussed clustering model"
$ws.Range("E46").Value2 = "- start looking at creating the  clustering model"
$ws.Range("D47").Value2 = "-EDA / Cleaning`n- Discussed clustering model"
$ws.Range("E47").Value2 = "-Finish resole cleaning and assit Rawad if required"
$ws.Range("D48").Value2 = "- Logged feedback for meeting with Phil`n- Discussed clustering model"
$ws.Range("E48").Value2 = "- Change up the report to include revised work tasks"
$ws.Range("D49").Value2 = "- setup Git repository , uploaded all documents so far , and shared with group`n-Submitted initial draft report"
$ws.Range("E49").Value2 = "- push meeting minutes to GIT `n- Review Mourads code `n-Assist others where needed"
